$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Exam A")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Exam A 2"
$ws.Cells.Clear()

$ws.Range("D1").Value = "Vishing"
$ws.Range("D2").Value = "Spoofing"
$ws.Range("D3").Value = "On-path"
$ws.Range("D4").Value = "DDoS"
$ws.Range("D5").Value = "Hoax"

$ws.Range("E1").Value = "Environmental Sensors"
$ws.Range("E2").Value = "Locking cabinets"
$ws.Range("E3").Value = "Video surveillance"
$ws.Range("E4").Value = "Full-disk encryption"
$ws.Range("E5").Value = "Biometric reader"
$ws.Range("E6").Value = "Smart card"
$ws.Range("E7").Value = "Cable lock"

$ws.Range("H1").Value = "HTTPS"
$ws.Range("H2").Value = "NTP"
$ws.Range("H3").Value = "SSH"
$ws.Range("H4").Value = "SRTP"
$ws.Range("H5").Value = "SNMPv3"

$ws.Range("I1").Value = "Something you have"
$ws.Range("I2").Value = "Something you know"
$ws.Range("I3").Value = "Something you can do"
$ws.Range("I4").Value = "Something you are"
$ws.Range("I5").Value = "Somewhere you are "

$ws.Range("L1").Value = "10.1.1.2 10.2.1.20 TCP 389 ALLOW"
$ws.Range("L2").Value = "10.2.1.33 10.1.1.7 TCP 443 ALLOW"
$ws.Range("L3").Value = "10.2.1.47 10.1.1.3 TCP 22 ALLOW"

$ws.Range("A6").Value = "B"
$ws.Range("A7").Value = "AC"
$ws.Range("A8").Value = "A"
$ws.Range("A9").Value = "AD"
$ws.Range("A10").Value = "C"
$ws.Range("A11").Value = "C"
$ws.Range("A12").Value = "A"
$ws.Range("A13").Value = "D"
$ws.Range("A14").Value = "AEG"
$ws.Range("A15").Value = "B"
$ws.Range("A16").Value = "C"
$ws.Range("A17").Value = "A"
$ws.Range("A18").Value = "A"
$ws.Range("A19").Value = "A"
$ws.Range("A20").Value = "C"
$ws.Range("A21").Value = "BE"
$ws.Range("A22").Value = "C"
$ws.Range("A23").Value = "BD"
$ws.Range("A24").Value = "B"
$ws.Range("A25").Value = "C"
$ws.Range("A26").Value = "B"
$ws.Range("A27").Value = "A"
$ws.Range("A28").Value = "B"
$ws.Range("A29").Value = "C"
$ws.Range("A30").Value = "F"
$ws.Range("A31").Value = "D"
$ws.Range("A32").Value = "AF"
$ws.Range("A33").Value = "B"
$ws.Range("A34").Value = "DF"
$ws.Range("A35").Value = "A"
$ws.Range("A36").Value = "D"
$ws.Range("A37").Value = "D"
$ws.Range("A38").Value = "DE"
$ws.Range("A39").Value = "D"
$ws.Range("A40").Value = "C"
$ws.Range("A41").Value = "C"
$ws.Range("A42").Value = "A"
$ws.Range("A43").Value = "B"
$ws.Range("A44").Value = "D"
$ws.Range("A45").Value = "C"
$ws.Range("A46").Value = "D"
$ws.Range("A47").Value = "A"
$ws.Range("A48").Value = "D"
$ws.Range("A49").Value = "C"
$ws.Range("A50").Value = "B"
$ws.Range("A51").Value = "D"
$ws.Range("A52").Value = "C"
$ws.Range("A53").Value = "B"
$ws.Range("A54").Value = "D"
$ws.Range("A55").Value = "C"
$ws.Range("A56").Value = "A"
$ws.Range("A57").Value = "B"
$ws.Range("A58").Value = "A"
$ws.Range("A59").Value = "A"
$ws.Range("A60").Value = "B"
$ws.Range("A61").Value = "B"
$ws.Range("A62").Value = "D"
$ws.Range("A63").Value = "B"
$ws.Range("A64").Value = "B"
$ws.Range("A65").Value = "C"
$ws.Range("A66").Value = "D"
$ws.Range("A67").Value = "A"
$ws.Range("A68").Value = "A"
$ws.Range("A69").Value = "C"
$ws.Range("A70").Value = "D"
$ws.Range("A71").Value = "D"
$ws.Range("A72").Value = "A"
$ws.Range("A73").Value = "D"
$ws.Range("A74").Value = "A"
$ws.Range("A75").Value = "D"
$ws.Range("A76").Value = "CD"
$ws.Range("A77").Value = "D"
$ws.Range("A78").Value = "C"
$ws.Range("A79").Value = "A"
$ws.Range("A80").Value = "BC"
$ws.Range("A81").Value = "B"
$ws.Range("A82").Value = "C"
$ws.Range("A83").Value = "C"
$ws.Range("A84").Value = "C"
$ws.Range("A85").Value = "B"
$ws.Range("A86").Value = "B"
$ws.Range("A87").Value = "A"
$ws.Range("A88").Value = "B"
$ws.Range("A89").Value = "A"
$ws.Range("A90").Value = "D"

$ws.Activate()
$ws.Range("L3").Select()